$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row ---
# Columns A:J were suffixed "_old" (the FV2404 / old-version side of the diff),
# columns L:U were suffixed "_new" (the FV2410 / new-version side of the diff).
# Column K just holds the literal header "diff" and stays untouched.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = [string]$cell.Value2
    if ($header -like "*_old") {
        $cell.Value = ($header.Substring(0, $header.Length - 4) + "_FV2404")
    } elseif ($header -like "*_new") {
        $cell.Value = ($header.Substring(0, $header.Length - 4) + "_FV2410")
    }
}

# --- Turn the used range (the whole sheet) into a table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U74"), $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
Write-Host "Applied header rename, table and frozen pane."
